# Update "Impor Data PD" sheet: fix row 2's running number, add a new
# data row (UDIN) as row 4, and extend the two list data validations
# down to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2, column A ("ID SEKOLAH" counter) changes from 3 to 1
$ws.Range("A2").Value = 1

# Build the new row 4 by cloning row 3 (same style set) then overwriting
# only the cells whose values actually differ for the new student record.
$ws.Range("A3:O3").Copy($ws.Range("A4:O4"))

$ws.Range("A4").Value = 1
$ws.Range("D4").Value = "'001456794"
$ws.Range("F4").Value = "UDIN"
$ws.Range("B4").Value = "'3205106707020009"
$ws.Range("E4").Value = 171814309

# Extend the two list-based data validations so they also cover row 4
$ws.Range("G2:G4").Validation.Delete()
$ws.Range("G2:G4").Validation.Add(3, 1, 1, '"Laki-laki,Perempuan"')

$ws.Range("J2:J4").Validation.Delete()
$ws.Range("J2:J4").Validation.Add(3, 1, 1, '"Islam,Kristen,Khatolik,Hindu,Budha"')

# Move the active selection to D7, matching the saved workbook state
$ws.Range("D7").Select()
